$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "James Trivette"
$ws.Range("B14").Value = "Federico Zoller | GREP"
$ws.Range("C14").Value = "Federico Manica | IMONTAGNA"
$ws.Range("D14").Value = "Alessandro Maffei | FC SAVIGNANO"
$ws.Range("E14").Value = "Federico Nicolodi | U.SGUARNA"
$ws.Range("F14").Value = "Alessio Debiasi | Mai una gioia"
